# Updates the cryptos worksheet cell values according to the latest
# crypto price/volume scrape, matching the target commit's changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value, whether the value must be
# forced to Text format so Excel does not silently reinterpret it as a
# number (which would corrupt formatted values like "1.00" or "0.993").
$changes = @(
    @("D2", "62.982.59", $false),
    @("E2", "  -2.08%  ", $false),
    @("D3", "2.676.03", $false),
    @("E3", "  -3.51%  ", $false),
    @("E4", "  -0.10%  ", $false),
    @("D5", "549.55", $true),
    @("E5", "  -4.94%  ", $false),
    @("D6", "157.28", $true),
    @("E6", "  -2.28%  ", $false),
    @("D7", "1.00", $true),
    @("E7", "  +0.28%  ", $false),
    @("D8", "0.585", $true),
    @("E8", "  -2.94%  ", $false),
    @("E9", "  -5.06%  ", $false),
    @("E10", "  -1.99%  ", $false),
    @("E11", "  -5.37%  ", $false),
    @("D12", "5.07", $true),
    @("E12", "  -14.12%  ", $false),
    @("D13", "3.148.29", $false),
    @("E13", "  -3.54%  ", $false),
    @("D14", "26.05", $true),
    @("E14", "  -4.65%  ", $false),
    @("D15", "62.856.24", $false),
    @("E15", "  -1.77%  ", $false),
    @("E16", "  -5.44%  ", $false),
    @("D17", "2.675.86", $false),
    @("E17", "  -4.04%  ", $false),
    @("D18", "11.92", $true),
    @("E18", "  -2.66%  ", $false),
    @("D19", "4.55", $true),
    @("E19", "  -6.63%  ", $false),
    @("E20", "  -5.73%  ", $false),
    @("D21", "6.28", $true),
    @("E21", "  -6.40%  ", $false),
    @("D22", "0.996", $true),
    @("E22", "  -0.45%  ", $false),
    @("D23", "0.504", $true),
    @("E23", "  -5.44%  ", $false),
    @("D24", "63.29", $true),
    @("E24", "  -3.00%  ", $false),
    @("E25", "  -2.18%  ", $false),
    @("D26", "0.993", $true),
    @("E26", "  -0.83%  ", $false),
    @("D27", "8.09", $true),
    @("E27", "  -6.49%  ", $false),
    @("D28", "0.0₃0848", $false),
    @("E28", "  -7.94%  ", $false),
    @("E29", "  -3.26%  ", $false),
    @("D30", "1.34", $true),
    @("E30", "  -2.32%  ", $false),
    @("D31", "6.98", $true),
    @("E31", "  -5.78%  ", $false),
    @("D32", "165.58", $true),
    @("E32", "  -0.73%  ", $false),
    @("D33", "0.999", $true),
    @("D34", "4.79", $true),
    @("E34", "  -4.77%  ", $false),
    @("E35", "  -4.02%  ", $false),
    @("D36", "1.43", $true),
    @("E36", "  -6.99%  ", $false),
    @("D37", "1.77", $true),
    @("E37", "  -3.25%  ", $false),
    @("D38", "339.63", $true),
    @("E38", "  -3.49%  ", $false),
    @("B39", "RenderToken", $false),
    @("C39", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", $false),
    @("D39", "6.15", $true),
    @("E39", "  -3.56%  ", $false),
    @("B40", "SuiNetwork", $false),
    @("C40", "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui", $false),
    @("D40", "0.933", $true),
    @("E40", "  -7.94%  ", $false),
    @("D41", "38.04", $true),
    @("E41", "  -2.87%  ", $false),
    @("D42", "3.92", $true),
    @("E42", "  -7.17%  ", $false),
    @("B43", "EnergySwap", $false),
    @("C43", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", $false),
    @("D43", "20.21", $true),
    @("E43", "  -6.99%  ", $false),
    @("B44", "InjectiveProtocol", $false),
    @("C44", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", $false),
    @("D44", "20.68", $true),
    @("E44", "  -9.23%  ", $false),
    @("B45", "Hedera", $false),
    @("C45", "https://coinranking.com/coin/jad286TjB+hedera-hbar", $false),
    @("D45", "0.0559", $true),
    @("E45", "  -6.52%  ", $false),
    @("B46", "Mantle", $false),
    @("C46", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", $false),
    @("D46", "0.613", $true),
    @("E46", "  -3.23%  ", $false),
    @("B47", "FirstDigitalUSD", $false),
    @("C47", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", $false),
    @("D47", "0.997", $true),
    @("E47", "  +0.00%  ", $false),
    @("E48", "  +0.00%  ", $false),
    @("E49", "  -4.35%  ", $false),
    @("D50", "128.71", $true),
    @("E50", "  -6.43%  ", $false),
    @("D51", "0.0239", $true),
    @("E51", "  -6.47%  ", $false)
)

foreach ($chg in $changes) {
    $ref = $chg[0]
    $val = $chg[1]
    $needsText = $chg[2]
    $cell = $ws.Range($ref)
    if ($needsText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $val
}
